# Update "想去人数" (interest count) values in column F on the
# "展览" and "全部类型" worksheets, matching the latest scraped data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 312
    3  = 94
    4  = 389
    5  = 11495
    6  = 747
    7  = 110
    10 = 146
    11 = 164
    14 = 48
    16 = 33
    17 = 323
    18 = 1304
    19 = 68
    20 = 896
    21 = 111
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
